$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows above the existing "salary" row (row 2), pushing it down to row 4.
# The inserted rows inherit formatting from row 2 (so column C keeps the date style).
$ws.Rows("2:3").Insert()

# Copy the date-formatted style that "salary"'s date cell (now C4) carries down to the
# two newly inserted date cells, so C2/C3 reuse the existing date style instead of
# a brand new number-format style getting created.
$ws.Range("C4").Copy()
$ws.Range("C2").PasteSpecial(-4122)
$ws.Range("C3").PasteSpecial(-4122)

# Row 2: Dividend from stocks
$ws.Range("A2").Value = "Dividend from stocks"
$ws.Range("B2").Value = 15000
$ws.Range("C2").Value = 45820.22928240741

# Row 3: Interest from fd
$ws.Range("A3").Value = "Interest from fd"
$ws.Range("B3").Value = 5000
$ws.Range("C3").Value = 45820.22928240741

# Row 4: salary (shifted down by the insert above; update its amount and date)
$ws.Range("B4").Value = 10000
$ws.Range("C4").Value = 45813.22928240741
